# Generate Report for Handoff
# Re-sorts the two tracked files (c09d6bce*, dd7c7456*) in each sheet and
# refreshes the c09d6bce row with the new "Ready for handoff" status produced
# by this handoff run.

$wb = $excel.ActiveWorkbook

$c09 = "c09d6bce-f51b-4dad-a143-c9492fc9e80a.md"
$dd7 = "dd7c7456-401e-4e27-989f-0d3c749511ba.md"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $dd7
$ws1.Range("B2").Value = "e2e\" + $dd7
$ws1.Range("A3").Value = $c09
$ws1.Range("B3").Value = "e2e\" + $c09

$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-22 11:37:21"

foreach ($h in $ws1.Hyperlinks) {
    $addr = [string]$h.Range().Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\" + $dd7
    }
    if ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\" + $c09
    }
}

$ws1.Columns.Item(5).ColumnWidth = 16.333333
$ws1.Columns.Item(6).ColumnWidth = 16.333333

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $dd7
$ws2.Range("G2").Value = "dd7c7456-401e-4e27-989f-0d3c749511ba.3fbe209e21cac527bf89dd7e897664ae027e3a20.zh-cn.xlf"

$ws2.Range("A3").Value = $c09
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "mt"
$ws2.Range("G3").Value = "c09d6bce-f51b-4dad-a143-c9492fc9e80a.4956508b5c9cc7ea61da3f91e1d41e32f9c1d930.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-22 11:37:16"

foreach ($h in $ws2.Hyperlinks) {
    $addr = [string]$h.Range().Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $dd7
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = $c09
    }
}

$ws2.Columns.Item(3).ColumnWidth = 16.333333

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $dd7
$ws3.Range("G2").Value = "dd7c7456-401e-4e27-989f-0d3c749511ba.3fbe209e21cac527bf89dd7e897664ae027e3a20.de-de.xlf"

$ws3.Range("A3").Value = $c09
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "mt"
$ws3.Range("G3").Value = "c09d6bce-f51b-4dad-a143-c9492fc9e80a.4956508b5c9cc7ea61da3f91e1d41e32f9c1d930.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-22 11:37:21"

foreach ($h in $ws3.Hyperlinks) {
    $addr = [string]$h.Range().Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $dd7
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = $c09
    }
}

$ws3.Columns.Item(3).ColumnWidth = 16.333333
